# Generate Report for Handback
# Update generated/handback timestamps (and priority flag) produced by a fresh
# handback-status report run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the 96a6b0b5.../c1525d83... rows
$wsOverview.Range("G4").Value = "2016-08-30 08:16:52"
$wsOverview.Range("G5").Value = "2016-08-30 08:16:52"

# zh-cn / de-de sheets: Priority column changes from "ht" to "mt"
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H4").Value = "2016-08-30 08:16:46"
$wsZhCn.Range("H5").Value = "2016-08-30 08:16:46"
$wsZhCn.Range("K4").Value = "2016-08-30 08:17:26"
$wsZhCn.Range("K5").Value = "2016-08-30 08:17:26"

# de-de sheet: Correspond Handoff Datetime (shares the same value as Overview!G4/G5)
$wsDeDe.Range("H4").Value = "2016-08-30 08:16:52"
$wsDeDe.Range("H5").Value = "2016-08-30 08:16:52"

# de-de sheet: Correspond Handback DateTime
$wsDeDe.Range("K4").Value = "2016-08-30 08:17:34"
$wsDeDe.Range("K5").Value = "2016-08-30 08:17:34"
